$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Pago" subscription field as column L
$ws.Range("L1").Value = "Pago"
$ws.Range("L2").Value = "Pendente"
$ws.Range("L3").Value = "Pendente"
$ws.Range("L4").Value = "Pendente"
$ws.Range("L5").Value = "Pendente"
$ws.Range("L6").Value = "Pendente"

# Update selection to reflect new column
$ws.Range("L:L").Select()
